$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.022.97'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -7.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.411.89'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -8.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '273.44'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.80%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3715'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -4.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3077'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.37'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -8.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.000'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06570'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -8.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.432'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.191'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '17.01'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -8.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.405.61'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -9.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001010'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -8.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05765'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -12.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.88'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -11.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.606'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -8.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.48'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.87'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.328'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '19.955.39'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -8.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.281'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '138.60'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.89'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.565.24'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -8.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.21'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -7.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.814'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -21.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.406'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -8.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8545'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -11.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07717'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.466'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05806'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.841'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.55%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1926'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02043'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.33'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.062'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -9.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.273'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -11.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5296'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -7.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.524'
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.21'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5135'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.803'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.09'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.049'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -9.58%  '
$ws.Range("E51").Value = '  +0.11%  '
